# Update scripts with new TPM values.
# The LR-pairs sheet previously had 2 data rows (one per Sending cluster: ECs, FAPs),
# each reporting only one Target cluster (MuSCs). The new TPM computation adds a
# second Target cluster (ECs) for each Sending cluster, expanding the sheet from
# 2 data rows to 4 data rows, and recomputing the specificity-derived columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Sending=ECs, Target=ECs (new target for this sending cluster)
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Rln3"
$ws.Range("C2").Value = "Rxfp1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.225097
$ws.Range("H2").Value = 0.675291
$ws.Range("I2").Value = 0.1920188148530651
$ws.Range("J2").Value = 0.1920188148530651
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.02528233333333333
$ws.Range("N2").Value = 0.075847
$ws.Range("O2").Value = 0.2893587312729617
$ws.Range("P2").Value = 0.2893587312729617
$ws.Range("Q2").Value = 0.005690977386333333
$ws.Range("R2").Value = 0.051218796477
$ws.Range("S2").Value = 0.05556232064642067
$ws.Range("T2").Value = 0.05556232064642066

# Row 3: Sending=ECs, Target=MuSCs (was row 2 before, target column re-specified)
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Rln3"
$ws.Range("C3").Value = "Rxfp1"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.225097
$ws.Range("H3").Value = 0.675291
$ws.Range("I3").Value = 0.1920188148530651
$ws.Range("J3").Value = 0.1920188148530651
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.06209133333333333
$ws.Range("N3").Value = 0.186274
$ws.Range("O3").Value = 0.7106412687270383
$ws.Range("P3").Value = 0.7106412687270383
$ws.Range("Q3").Value = 0.01397657285933333
$ws.Range("R3").Value = 0.125789155734
$ws.Range("S3").Value = 0.1364564942066445
$ws.Range("T3").Value = 0.1364564942066444

# Row 4: Sending=FAPs, Target=ECs (new row)
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Rln3"
$ws.Range("C4").Value = "Rxfp1"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.9471683333333334
$ws.Range("H4").Value = 2.841505
$ws.Range("I4").Value = 0.807981185146935
$ws.Range("J4").Value = 0.8079811851469348
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02528233333333333
$ws.Range("N4").Value = 0.075847
$ws.Range("O4").Value = 0.2893587312729617
$ws.Range("P4").Value = 0.2893587312729617
$ws.Range("Q4").Value = 0.02394662552611111
$ws.Range("R4").Value = 0.215519629735
$ws.Range("S4").Value = 0.2337964106265411
$ws.Range("T4").Value = 0.2337964106265411

# Row 5: Sending=FAPs, Target=MuSCs (new row, holds the old row 3 values)
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Rln3"
$ws.Range("C5").Value = "Rxfp1"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.9471683333333334
$ws.Range("H5").Value = 2.841505
$ws.Range("I5").Value = 0.807981185146935
$ws.Range("J5").Value = 0.8079811851469348
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.06209133333333333
$ws.Range("N5").Value = 0.186274
$ws.Range("O5").Value = 0.7106412687270383
$ws.Range("P5").Value = 0.7106412687270383
$ws.Range("Q5").Value = 0.05881094470777778
$ws.Range("R5").Value = 0.52929850237
$ws.Range("S5").Value = 0.5741847745203938
$ws.Range("T5").Value = 0.5741847745203937
